$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = $false
$ws.Range("B2").Value = 'commonswiki'
$ws.Range("C2").Value = 'Anvilaquarius'
$ws.Range("D2").Value = "'1661"
$ws.Range("C3").Value = 'QuickStatementsBot'
$ws.Range("D3").Value = "'974"
$ws.Range("C4").Value = 'Mr.Ibrahembot'
$ws.Range("D4").Value = "'665"
$ws.Range("A5").Value = $false
$ws.Range("C5").Value = 'Vojtěch Dostál'
$ws.Range("D5").Value = "'653"
$ws.Range("C6").Value = 'Research Bot'
$ws.Range("D6").Value = "'638"
$ws.Range("B7").Value = 'commonswiki'
$ws.Range("C7").Value = 'SchlurcherBot'
$ws.Range("D7").Value = "'286"
$ws.Range("B8").Value = 'commonswiki'
$ws.Range("C8").Value = 'KolbertBot'
$ws.Range("D8").Value = "'285"
$ws.Range("B9").Value = 'cebwiki'
$ws.Range("C9").Value = 'Lsjbot'
$ws.Range("D9").Value = "'253"
$ws.Range("B10").Value = 'arwiki'
$ws.Range("C10").Value = 'JarBot'
$ws.Range("D10").Value = "'233"
$ws.Range("A11").Value = $true
$ws.Range("B11").Value = 'wikidatawiki'
$ws.Range("C11").Value = 'Dcirovicbot'
$ws.Range("D11").Value = "'181"
$ws.Range("B12").Value = 'wikidatawiki'
$ws.Range("C12").Value = 'Edoderoobot'
$ws.Range("D12").Value = "'130"
$ws.Range("A13").Value = $true
$ws.Range("C13").Value = 'Geograph Update Bot'
$ws.Range("D13").Value = "'115"
$ws.Range("B14").Value = 'dewiki'
$ws.Range("C14").Value = 'Drüfft'
$ws.Range("D14").Value = "'67"
$ws.Range("C15").Value = 'Arnaud Palastowicz'
$ws.Range("D15").Value = "'58"
$ws.Range("B16").Value = 'commonswiki'
$ws.Range("C16").Value = 'Reykholt'
$ws.Range("D16").Value = "'54"
$ws.Range("B17").Value = 'commonswiki'
$ws.Range("C17").Value = 'Martin H.'
$ws.Range("D17").Value = "'52"
$ws.Range("B18").Value = 'commonswiki'
$ws.Range("C18").Value = 'INS Pirat'
$ws.Range("D18").Value = "'41"
$ws.Range("B19").Value = 'plwiktionary'
$ws.Range("C19").Value = '84.10.92.69'
$ws.Range("D19").Value = "'37"
$ws.Range("B20").Value = 'enwiki'
$ws.Range("C20").Value = 'Geekkurosaki'
$ws.Range("D20").Value = "'36"
$ws.Range("B21").Value = 'commonswiki'
$ws.Range("C21").Value = 'Mr.Nostalgic'
$ws.Range("D21").Value = "'34"
